# Update the NATMI ligand-receptor TPM-derived metrics on Sheet1 (rows 2-6)
# with newly recomputed values ("update scripts wuth new tpm").
#
# Columns touched (1-based / letter):
#   G  Ligand average expression value
#   H  Ligand total expression value
#   I  Ligand derived specificity of average expression value
#   J  Ligand derived specificity of total expression value
#   M  Receptor average expression value
#   Q  Edge average expression weight
#   R  Edge total expression weight
#   S  Edge average expression derived specificity
#   T  Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> C3/Cd19)
$ws.Range("G2").Value = 3.820425
$ws.Range("H2").Value = 11.461275
$ws.Range("I2").Value = 0.02049663039797357
$ws.Range("J2").Value = 0.02049663039797357
$ws.Range("M2").Value = 0.04150766666666666
$ws.Range("Q2").Value = 0.158576927425
$ws.Range("R2").Value = 1.427192346825
$ws.Range("S2").Value = 0.02049663039797357
$ws.Range("T2").Value = 0.02049663039797357

# Row 3 (FAPs -> C3/Cd19)
$ws.Range("I3").Value = 0.7542622677884155
$ws.Range("J3").Value = 0.7542622677884157
$ws.Range("M3").Value = 0.04150766666666666
$ws.Range("Q3").Value = 5.835524697285111
$ws.Range("S3").Value = 0.7542622677884155
$ws.Range("T3").Value = 0.7542622677884157

# Row 4 (Inflammatory-Mac -> C3/Cd19)
$ws.Range("G4").Value = 30.51067
$ws.Range("H4").Value = 91.53201
$ws.Range("I4").Value = 0.1636901460399144
$ws.Range("J4").Value = 0.1636901460399144
$ws.Range("M4").Value = 0.04150766666666666
$ws.Range("Q4").Value = 1.266426720136667
$ws.Range("R4").Value = 11.39784048123
$ws.Range("S4").Value = 0.1636901460399144
$ws.Range("T4").Value = 0.1636901460399144

# Row 5 (MuSCs -> C3/Cd19)
$ws.Range("G5").Value = 0.258813
$ws.Range("H5").Value = 0.776439
$ws.Range("I5").Value = 0.001388535150720334
$ws.Range("J5").Value = 0.001388535150720334
$ws.Range("M5").Value = 0.04150766666666666
$ws.Range("Q5").Value = 0.010742723733
$ws.Range("R5").Value = 0.09668451359699999
$ws.Range("S5").Value = 0.001388535150720334
$ws.Range("T5").Value = 0.001388535150720334

# Row 6 (Resolving-Mac -> C3/Cd19)
$ws.Range("G6").Value = 11.213844
$ws.Range("H6").Value = 33.641532
$ws.Range("I6").Value = 0.0601624206229761
$ws.Range("J6").Value = 0.0601624206229761
$ws.Range("M6").Value = 0.04150766666666666
$ws.Range("Q6").Value = 0.465460498804
$ws.Range("R6").Value = 4.189144489236
$ws.Range("S6").Value = 0.0601624206229761
$ws.Range("T6").Value = 0.0601624206229761
